# Generate Report for Handback
# Update the localization status report: the zh-cn and de-de handback
# files are now in sync with en-US, so update their status, handback
# timestamps, and clear the stale "version mismatch" error detail.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns for both rows
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn detail sheet
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("L2").Value = "2017-02-17 09:59:08"
$wsZhCn.Range("L3").Value = "2017-02-17 09:59:08"
$wsZhCn.Range("R2").Value = ""

# de-de detail sheet
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("L2").Value = "2017-02-17 09:59:33"
$wsDeDe.Range("L3").Value = "2017-02-17 09:59:33"
$wsDeDe.Range("R2").Value = ""
